$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1868.1818
$ws.Range("I17").Value = 1500
$ws.Range("K17").Value = 4500
$ws.Range("M17").Value = -4332

$ws.Range("H18").Value = 18650
$ws.Range("I18").Value = 18650
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 18650
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -18366

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H86").Value = 4390451.5
$ws.Range("I86").Value = 4569
$ws.Range("J86").Value = 7523225
$ws.Range("K86").Value = 4569
$ws.Range("L86").Value = 7523225
$ws.Range("M86").Value = -3446
$ws.Range("N86").Value = -7525471

$ws.Range("H89").Value = 4390451.5
$ws.Range("I89").Value = 4569
$ws.Range("J89").Value = 7523225
$ws.Range("K89").Value = 22845
$ws.Range("L89").Value = 37616125
$ws.Range("M89").Value = -17229
$ws.Range("N89").Value = -37627357

$ws.Range("H106").Value = 4189.5
$ws.Range("I106").Value = 4835
$ws.Range("J106").Value = 2253
$ws.Range("K106").Value = 4835
$ws.Range("L106").Value = 2253
$ws.Range("M106").Value = -4204
$ws.Range("N106").Value = -3515

$ws.Range("H113").Value = 8394.385
$ws.Range("I113").Value = 5948.5
$ws.Range("J113").Value = 9481.444
$ws.Range("K113").Value = 5948.5
$ws.Range("L113").Value = 9481.444
$ws.Range("M113").Value = -2694.5
$ws.Range("N113").Value = -15989.444

$ws.Range("H125").Value = 9747.5
$ws.Range("I125").Value = 1229.1666
$ws.Range("J125").Value = 18265.834
$ws.Range("K125").Value = 11062.4994
$ws.Range("L125").Value = 164392.506
$ws.Range("M125").Value = -8602.499400000001
$ws.Range("N125").Value = -169312.506

$ws.Range("H138").Value = 3009.1333
$ws.Range("J138").Value = 5599.6
$ws.Range("L138").Value = 16798.8
$ws.Range("N138").Value = -27078.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18874790
$ws.Range("I32").Value = 18874790
$ws.Range("K32").Value = 18874790
$ws.Range("M32").Value = -18874503

$ws.Range("H45").Value = 3182
$ws.Range("I45").Value = 1937.3334
$ws.Range("J45").Value = 5671.3335
$ws.Range("K45").Value = 1937.3334
$ws.Range("L45").Value = 5671.3335
$ws.Range("M45").Value = -1560.3334
$ws.Range("N45").Value = -6425.3335

$ws.Range("H61").Value = 4024
$ws.Range("I61").Value = 3245.3125
$ws.Range("K61").Value = 3245.3125
$ws.Range("M61").Value = -3033.3125

$ws.Range("H102").Value = 30305422
$ws.Range("I102").Value = 2963.7144
$ws.Range("J102").Value = 83334720
$ws.Range("K102").Value = 2963.7144
$ws.Range("L102").Value = 83334720
$ws.Range("M102").Value = -1341.7144
$ws.Range("N102").Value = -83337964

$ws.Range("H122").Value = 4846.5625
$ws.Range("I122").Value = 3443.25
$ws.Range("J122").Value = 6249.875
$ws.Range("K122").Value = 10329.75
$ws.Range("L122").Value = 18749.625
$ws.Range("M122").Value = -7879.75
$ws.Range("N122").Value = -23649.625

$ws.Range("H136").Value = 4024
$ws.Range("I136").Value = 3245.3125
$ws.Range("K136").Value = 9735.9375
$ws.Range("M136").Value = -7185.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2492.0195
$ws.Range("I31").Value = 1640.5
$ws.Range("J31").Value = 4981.077
$ws.Range("K31").Value = 1640.5
$ws.Range("L31").Value = 4981.077
$ws.Range("M31").Value = -1345.5
$ws.Range("N31").Value = -5571.077

$ws.Range("H34").Value = 2492.0195
$ws.Range("I34").Value = 1640.5
$ws.Range("J34").Value = 4981.077
$ws.Range("K34").Value = 1640.5
$ws.Range("L34").Value = 4981.077
$ws.Range("M34").Value = -1438.5
$ws.Range("N34").Value = -5385.077

$ws.Range("H94").Value = 6369.8945
$ws.Range("I94").Value = 17303.5
$ws.Range("K94").Value = 17303.5
$ws.Range("M94").Value = -16852.5

$ws.Range("H122").Value = 396371.8
$ws.Range("I122").Value = 785921.3
$ws.Range("K122").Value = 2357763.9
$ws.Range("M122").Value = -2355313.9

$ws.Range("H134").Value = 5225.9023
$ws.Range("I134").Value = 3531.4285
$ws.Range("K134").Value = 10594.2855
$ws.Range("M134").Value = -8059.2855

$ws.Range("H141").Value = 28937.375
$ws.Range("J141").Value = 28937.375
$ws.Range("L141").Value = 28937.375
$ws.Range("N141").Value = -39297.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4135.143
$ws.Range("I126").Value = 1100.25
$ws.Range("J126").Value = 5349.1
$ws.Range("K126").Value = 3300.75
$ws.Range("L126").Value = 16047.3
$ws.Range("M126").Value = -830.75
$ws.Range("N126").Value = -20987.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2405
$ws.Range("I93").Value = 2546.2632
$ws.Range("K93").Value = 2546.2632
$ws.Range("M93").Value = -1298.2632

$ws.Range("H122").Value = 5444.7144
$ws.Range("I122").Value = 2720.7778
$ws.Range("J122").Value = 6735
$ws.Range("K122").Value = 8162.3334
$ws.Range("L122").Value = 20205
$ws.Range("M122").Value = -5712.3334
$ws.Range("N122").Value = -25105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 25123.75
$ws.Range("J51").Value = 40495
$ws.Range("L51").Value = 40495
$ws.Range("N51").Value = -41515

$ws.Range("H81").Value = 5850066.5
$ws.Range("I81").Value = 2348.5833
$ws.Range("J81").Value = 15874726
$ws.Range("K81").Value = 4697.1666
$ws.Range("L81").Value = 31749452
$ws.Range("M81").Value = -3636.1666
$ws.Range("N81").Value = -31751574

$ws.Range("H84").Value = 5850066.5
$ws.Range("I84").Value = 2348.5833
$ws.Range("J84").Value = 15874726
$ws.Range("K84").Value = 23485.833
$ws.Range("L84").Value = 158747260
$ws.Range("M84").Value = -18181.833
$ws.Range("N84").Value = -158757868

$ws.Range("H107").Value = 579.2778
$ws.Range("I107").Value = 401.6875
$ws.Range("K107").Value = 1205.0625
$ws.Range("M107").Value = 714.9375

$ws.Range("H113").Value = 336.69232
$ws.Range("I113").Value = 310.7143
$ws.Range("K113").Value = 932.1428999999999
$ws.Range("M113").Value = 1237.8571

$ws.Range("I122").Value = 529979.4
$ws.Range("J122").Value = 5486.615
$ws.Range("K122").Value = 1589938.2
$ws.Range("L122").Value = 16459.845
$ws.Range("M122").Value = -1587488.2
$ws.Range("N122").Value = -21359.845
